# Apply "add syllabus and textbook" edits to the ADS schedule workbook.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("schedule")
$wsSheet1   = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Sheet "schedule": weeks 10-14 (rows 12-16) get new readings / links, new
# homework assignments, and week 12's topic switches from generic "Parallel
# processing, batch jobs" to "Working on a computing cluster" with a new
# prepare item.
# ---------------------------------------------------------------------------

# Week 10 (row 12) - Tree based methods
$wsSchedule.Range("E12").Value = "* Read WOMD: Ch 3-4 (MB)`n* [When Data Science Destabilizes Democracy and Facilitates Genocide](http://www.fast.ai/2017/11/02/ethics/)"
$wsSchedule.Range("I12").Value = "* ML homework"

# Week 11 (row 13) - Big Data
$wsSchedule.Range("E13").Value = "* Read WOMD: Ch 5-6 (AS)`n* [More than a Million Pro-Repeal Net Neutrality Comments were Likely Faked](https://hackernoon.com/more-than-a-million-pro-repeal-net-neutrality-comments-were-likely-faked-e9f0e3ed36a6)"
$wsSchedule.Range("I13").Value = "* Big data homework"

# Week 12 (row 14) - Working on a computing cluster (was "Parallel processing, batch jobs")
$wsSchedule.Range("C14").Value = "Working on a computing cluster"
$wsSchedule.Range("D14").Value = "* Write a script to send to a remote server to perform an analysis. `n* Write a script that utilizes more than one computing core. "
$wsSchedule.Range("E14").Value = "* Read WOMD: Ch 7-8 (RA)`n* [Activision Patents Matchmaking That Encourages Players To Buy Microtransactions](https://kotaku.com/activision-patents-matchmaking-that-encourages-players-1819630937)"
$wsSchedule.Range("I14").Value = "* cluster computing homework"

# Week 13 (row 15) - Other Tools
$wsSchedule.Range("E15").Value = "* Read WOMD: Ch 9-10 (RAD)`n* Compare [[this article on practical tips for success with ML]](https://www.datanami.com/2018/01/17/practical-tips-success-machine-learning/) with this post on [[google automating ML]](https://www.datanami.com/2018/01/17/google-automate-machine-learning-automl-service/)"

# Week 14 (row 16) - Finals prep placeholder (Prepare column gets a single space)
$wsSchedule.Range("E16").Value = " "

# Row heights grow to fit the newly-wrapped text in these rows.
$wsSchedule.Rows.Item(12).RowHeight = 64.5
$wsSchedule.Rows.Item(13).RowHeight = 90
$wsSchedule.Rows.Item(14).RowHeight = 77.25
$wsSchedule.Rows.Item(15).RowHeight = 128.25

# ---------------------------------------------------------------------------
# Sheet "Sheet1": the points table gets three new Assignment/Learning pairs
# (ML, Big data, cluster computing) inserted ahead of the Blog Posts / Final
# Project rows, which shift down from rows 17-18 to rows 23-24.
# ---------------------------------------------------------------------------

# Preserve the two rows that get pushed down before overwriting them.
$blogPostsCategory = $wsSheet1.Range("C17").Value2
$finalProjectCategory = $wsSheet1.Range("C18").Value2

$wsSheet1.Range("A17").Value = 10
$wsSheet1.Range("B17").Value = "ML"
$wsSheet1.Range("C17").Value = "Assignment"
$wsSheet1.Range("D17").Value = 10

$wsSheet1.Range("A18").Value = 10
$wsSheet1.Range("B18").Value = "PR ML"
$wsSheet1.Range("C18").Value = "Learning"
$wsSheet1.Range("D18").Value = 5

$wsSheet1.Range("A19").Value = 11
$wsSheet1.Range("B19").Value = "Big data"
$wsSheet1.Range("C19").Value = "Assignment"
$wsSheet1.Range("D19").Value = 10

$wsSheet1.Range("A20").Value = 11
$wsSheet1.Range("B20").Value = "PR big data"
$wsSheet1.Range("C20").Value = "Learning"
$wsSheet1.Range("D20").Value = 5

$wsSheet1.Range("A21").Value = 12
$wsSheet1.Range("B21").Value = "cluster computing"
$wsSheet1.Range("C21").Value = "Assignment"
$wsSheet1.Range("D21").Value = 10

$wsSheet1.Range("A22").Value = 12
$wsSheet1.Range("B22").Value = "PR Cluster computing"
$wsSheet1.Range("C22").Value = "Learning"
$wsSheet1.Range("D22").Value = 5

$wsSheet1.Range("A23").Value = 15
$wsSheet1.Range("B23").Value = "Blog Posts"
$wsSheet1.Range("C23").Value = "Discussion"
$wsSheet1.Range("D23").Value = 50

$wsSheet1.Range("A24").Value = 16
$wsSheet1.Range("B24").Value = "Final Project"
$wsSheet1.Range("C24").Value = "Project"
$wsSheet1.Range("D24").Value = 50

# ---------------------------------------------------------------------------
# View state: scroll/selection changes recorded in the source diff.
# ---------------------------------------------------------------------------

$wsSheet1.Activate()
$wsSheet1.Range("H20").Select()

$wsSchedule.Activate()
$excel.ActiveWindow.ScrollRow = 8
$wsSchedule.Range("H11").Select()
